# Insert a new data row at row 586 in the "Femacal de La Calera - Ají"
# sheet, pushing the existing rows 586-609 down to 587-610, and fill in
# the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 586..609 down to 587..610 by inserting a blank row at 586.
$ws.Rows.Item(586).Insert()

# Populate the newly inserted row 586 with the new observation.
$ws.Range("A586").Value = 3
$ws.Range("B586").Value = "Femacal de La Calera"
$ws.Range("C586").Value = "Coquimbo"
$ws.Range("D586").Value = 44939
$ws.Range("E586").Value = 5
$ws.Range("F586").Value = 100112021
$ws.Range("G586").Value = "Ají"
$ws.Range("H586").Value = "Inferno"
$ws.Range("I586").Value = "Primera"
$ws.Range("J586").Value = 73
$ws.Range("K586").Value = 19000
$ws.Range("L586").Value = 19500
$ws.Range("M586").Value = 19260
$ws.Range("N586").Value = '$/caja 15 kilos'
$ws.Range("O586").Value = "Limache"
$ws.Range("P586").Value = 1284
$ws.Range("Q586").Value = 15
$ws.Range("R586").Value = "Hortaliza"
